$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new column before EE, shifting EE:FI (and beyond) one column to the right.
$ws.Columns("EE").Insert()

# Fill the new column: header (row 1) continues the December date sequence,
# data rows (2-25) get the same "-" placeholder used for missing data.
$ws.Range("EE1").Value = "02-dec"
$ws.Range("EE2:EE25").Value = "-"
